# Update the 4COM02_ProteomicsDataProcessing template's ontology terms.
#
# For each software name already listed in the "acquisition software" (B),
# "analysis software" (E) and "data processing software" (H) columns, fill
# in the matching Term Source REF ("MS") and Term Accession Number
# (a purl.obolibrary.org MS_* URL) columns that sit next to them.
# Also bump the template version in the metadata sheet from 1.1.4 to 1.1.5.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("4COM02_ProteomicsDataProcessing")
$ws2 = $wb.Worksheets.Item("SwateTemplateMetadata")

# Term Accession Number values for column D - "Parameter [acquisition software]" (rows 2-7)
$acquisitionTerms = @(
    "http://purl.obolibrary.org/obo/MS_1001483",
    "http://purl.obolibrary.org/obo/MS_1000688",
    "http://purl.obolibrary.org/obo/MS_1000551",
    "http://purl.obolibrary.org/obo/MS_1001877",
    "http://purl.obolibrary.org/obo/MS_1000532",
    "http://purl.obolibrary.org/obo/MS_1000706"
)

# Term Accession Number values for column G - "Parameter [analysis software]" (rows 2-7)
$analysisTerms = @(
    "http://purl.obolibrary.org/obo/MS_1000551",
    "http://purl.obolibrary.org/obo/MS_1001799",
    "http://purl.obolibrary.org/obo/MS_1002251",
    "http://purl.obolibrary.org/obo/MS_1000712",
    "http://purl.obolibrary.org/obo/MS_1000667",
    "http://purl.obolibrary.org/obo/MS_1000707"
)

# Term Accession Number values for column J - "Parameter [data processing software]" (rows 2-7)
$dataProcessingTerms = @(
    "http://purl.obolibrary.org/obo/MS_1000551",
    "http://purl.obolibrary.org/obo/MS_1000667",
    "http://purl.obolibrary.org/obo/MS_1000707",
    "http://purl.obolibrary.org/obo/MS_1000533",
    "http://purl.obolibrary.org/obo/MS_1002871",
    "http://purl.obolibrary.org/obo/MS_1000536"
)

# Fill Term Source REF (col C) + Term Accession Number (col D) for acquisition software
for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 2
    $ws1.Cells.Item($r, 3).Value = "MS"
    $ws1.Cells.Item($r, 4).Value = $acquisitionTerms[$i]
}

# Fill Term Source REF (col F) + Term Accession Number (col G) for analysis software
for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 2
    $ws1.Cells.Item($r, 6).Value = "MS"
    $ws1.Cells.Item($r, 7).Value = $analysisTerms[$i]
}

# Fill Term Source REF (col I) + Term Accession Number (col J) for data processing software
for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 2
    $ws1.Cells.Item($r, 9).Value = "MS"
    $ws1.Cells.Item($r, 10).Value = $dataProcessingTerms[$i]
}

# Column layout: widen column B (now bestFit content is longer), and split
# the previously-merged, zero-width hidden columns I/J (9-10) into their own
# widths now that they hold real (hidden helper) content.
$ws1.Columns.Item(2).ColumnWidth = 38.8
$ws1.Columns.Item(9).ColumnWidth = 29.6
$ws1.Columns.Item(10).ColumnWidth = 36.6

# Bump template version 1.1.4 -> 1.1.5 on the metadata sheet (leading
# apostrophe keeps it stored/typed as text, matching the original cell).
$ws2.Range("B3").Value = "'1.1.5"
